$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "marginalized particle filter: anchoring can match long-term inflation forecasts w/o targeting them (can explain Missing Deflation/Inflation)"
$ws.Range("B4").Value = "Bayesian, constant gain I think: can match intertia in inflation"
$ws.Range("B5").Value = "Bayesian, endogenous gain, model matches time-varying volatility of macro variables, high in Great Inflation, falling in Great Moderation"

$ws.Range("A9").Value = "Eusepi Preston 2011 AER"
$ws.Range("B9").Value = "outlier b/c calibrated, but their whole point is the quantitative performance is better than RE version, e.g. persistence and hum-shaped IRFs despite iid shocks."
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

$ws.Range("B10").Select()
